# Apply edits described by the commit: fix surplus number (K column factor
# 1.0565 -> 1) on the "Customer Quote" sheet, and update the active cell
# selection to G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

# Fix surplus number: change the 1.0565 factor to 1 for rows 16, 17, 22, 25, 28, 29, 33
$rowsToFix = @(16, 17, 22, 25, 28, 29, 33)
foreach ($row in $rowsToFix) {
    $cell = $ws.Range("K$row")
    $cell.Value = 1
}

# Update the selected/active cell on the sheet (was A33, now G6)
$ws.Activate()
$ws.Range("G6").Select()

$wb.Save()
